# Auto-generated Excel COM-interop script
# Applies "Actualizacion automatica 2025-06-01 08:00:06" changes to
# sheets "VENTAS POR GRUPO" and "VENTA MENSUAL"

$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
# All non-zero monthly sales figures (rows 2-54) are reset to 0,
# and the "N de 53" summary labels in row 55 follow suit.
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("K5").Value = 0
$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = 0
$ws1.Range("K6").Value = 0
$ws1.Range("L6").Value = 0
$ws1.Range("K14").Value = 0
$ws1.Range("L14").Value = 0
$ws1.Range("N14").Value = 0
$ws1.Range("C22").Value = 0
$ws1.Range("L22").Value = 0
$ws1.Range("L24").Value = 0
$ws1.Range("C26").Value = 0
$ws1.Range("D26").Value = 0
$ws1.Range("K26").Value = 0
$ws1.Range("L26").Value = 0
$ws1.Range("D27").Value = 0
$ws1.Range("L27").Value = 0
$ws1.Range("D28").Value = 0
$ws1.Range("L28").Value = 0
$ws1.Range("E29").Value = 0
$ws1.Range("K29").Value = 0
$ws1.Range("L29").Value = 0
$ws1.Range("D36").Value = 0
$ws1.Range("K36").Value = 0
$ws1.Range("L36").Value = 0
$ws1.Range("M36").Value = 0
$ws1.Range("D39").Value = 0
$ws1.Range("L41").Value = 0
$ws1.Range("C43").Value = 0
$ws1.Range("D43").Value = 0
$ws1.Range("L43").Value = 0
$ws1.Range("M43").Value = 0
$ws1.Range("L44").Value = 0
$ws1.Range("L45").Value = 0
$ws1.Range("G47").Value = 0
$ws1.Range("D50").Value = 0
$ws1.Range("L50").Value = 0
$ws1.Range("N51").Value = 0
$ws1.Range("C54").Value = 0
$ws1.Range("C55").Value = "0 de 53"
$ws1.Range("D55").Value = "0 de 53"
$ws1.Range("E55").Value = "0 de 53"
$ws1.Range("G55").Value = "0 de 53"
$ws1.Range("K55").Value = "0 de 53"
$ws1.Range("L55").Value = "0 de 53"
$ws1.Range("M55").Value = "0 de 53"
$ws1.Range("N55").Value = "0 de 53"

# --- Sheet 2: VENTA MENSUAL ---
# Monthly columns roll forward one month (febrero..mayo -> marzo..junio);
# column headers shift and each row's C/D/E values move left one column
# while column F (the new month) starts at 0.
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("C1").Value = "marzo"
$ws2.Range("D1").Value = "abril"
$ws2.Range("E1").Value = "mayo"
$ws2.Range("F1").Value = "junio"
$ws2.Range("C4").Value = 0
$ws2.Range("C5").Value = 12803.66
$ws2.Range("D5").Value = 3471.96
$ws2.Range("E5").Value = 4158.27
$ws2.Range("F5").Value = 0
$ws2.Range("C6").Value = 778.48
$ws2.Range("D6").Value = 710.14
$ws2.Range("E6").Value = 1528.39
$ws2.Range("F6").Value = 0
$ws2.Range("C11").Value = 0
$ws2.Range("D11").Value = 262.99
$ws2.Range("E11").Value = 0
$ws2.Range("E14").Value = 2169.75
$ws2.Range("F14").Value = 0
$ws2.Range("C18").Value = 0
$ws2.Range("C19").Value = 4084.01
$ws2.Range("D19").Value = 6506.43
$ws2.Range("E19").Value = 0
$ws2.Range("C22").Value = 7471.12
$ws2.Range("D22").Value = 607.78
$ws2.Range("E22").Value = 4141.42
$ws2.Range("F22").Value = 0
$ws2.Range("C24").Value = 890.88
$ws2.Range("D24").Value = 1567.1
$ws2.Range("E24").Value = 61.78
$ws2.Range("F24").Value = 0
$ws2.Range("D25").Value = 129.6
$ws2.Range("E25").Value = 0
$ws2.Range("C26").Value = 13296.81
$ws2.Range("D26").Value = 29332.26
$ws2.Range("E26").Value = 6249.76
$ws2.Range("F26").Value = 0
$ws2.Range("C27").Value = 14695.45
$ws2.Range("D27").Value = 7533.56
$ws2.Range("E27").Value = 7315.29
$ws2.Range("F27").Value = 0
$ws2.Range("C28").Value = 4901.56
$ws2.Range("D28").Value = 5978.87
$ws2.Range("E28").Value = 3897.23
$ws2.Range("F28").Value = 0
$ws2.Range("C29").Value = 9655.45
$ws2.Range("D29").Value = 7135.59
$ws2.Range("E29").Value = 3563.29
$ws2.Range("F29").Value = 0
$ws2.Range("C31").Value = 0
$ws2.Range("D31").Value = 4280.48
$ws2.Range("E31").Value = 0
$ws2.Range("C35").Value = 0
$ws2.Range("C36").Value = 615.17
$ws2.Range("D36").Value = 722.57
$ws2.Range("E36").Value = 5704.92
$ws2.Range("F36").Value = 0
$ws2.Range("C39").Value = 1961.95
$ws2.Range("D39").Value = 2673.89
$ws2.Range("E39").Value = 2403.41
$ws2.Range("F39").Value = 0
$ws2.Range("C41").Value = 2526.46
$ws2.Range("D41").Value = 0
$ws2.Range("E41").Value = 2689.09
$ws2.Range("F41").Value = 0
$ws2.Range("C42").Value = 860.17
$ws2.Range("D42").Value = 15577.98
$ws2.Range("E42").Value = 0
$ws2.Range("C43").Value = 3687.17
$ws2.Range("D43").Value = 3716.29
$ws2.Range("E43").Value = 7574.03
$ws2.Range("F43").Value = 0
$ws2.Range("C44").Value = 814.08
$ws2.Range("D44").Value = 418.61
$ws2.Range("E44").Value = 731.63
$ws2.Range("F44").Value = 0
$ws2.Range("D45").Value = 1091.58
$ws2.Range("E45").Value = 722.54
$ws2.Range("F45").Value = 0
$ws2.Range("C47").Value = 4821.54
$ws2.Range("D47").Value = 969.61
$ws2.Range("E47").Value = 798
$ws2.Range("F47").Value = 0
$ws2.Range("C49").Value = 3617.78
$ws2.Range("D49").Value = 0
$ws2.Range("C50").Value = 4748.11
$ws2.Range("D50").Value = 4141.1
$ws2.Range("E50").Value = 4953.13
$ws2.Range("F50").Value = 0
$ws2.Range("D51").Value = 3336.39
$ws2.Range("E51").Value = -11.75
$ws2.Range("F51").Value = 0
$ws2.Range("E54").Value = 144
$ws2.Range("F54").Value = 0
$ws2.Range("C55").Value = 92229.85
$ws2.Range("D55").Value = 100164.78
$ws2.Range("E55").Value = 58794.18
$ws2.Range("F55").Value = 0

# Column widths on VENTA MENSUAL were adjusted to fit the new headers.
# ColumnWidth uses "characters" units that are 5/6 narrower than the
# stored worksheet width, so we compensate when setting explicit widths.
$colWidthOffset = 5 / 6
$ws2.Columns.Item(3).ColumnWidth = 14 - $colWidthOffset
$ws2.Columns.Item(4).ColumnWidth = 14 - $colWidthOffset
$ws2.Columns.Item(5).ColumnWidth = 13 - $colWidthOffset
$ws2.Columns.Item(6).ColumnWidth = 11 - $colWidthOffset

